$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header labels (German -> English)
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Years"
$ws.Range("D1").Value = "Town"

Write-Output "done"
